# Scout_Config.xlsx update: append new colormap options to the
# Cmap_Options list on Sheet1 (A6:A9): Reds, Blues, Greens, Oranges.
#
# Shared-string insertion order matters for a byte-faithful round trip,
# so "Blues" (A7) is written before "Reds" (A6) to match the order the
# strings were originally introduced into xl/sharedStrings.xml.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Blues"
$ws.Range("A6").Value = "Reds"
$ws.Range("A8").Value = "Greens"
$ws.Range("A9").Value = "Oranges"

# Move the active selection down to A10, just past the newly written data,
# matching where the cursor ended up after the edit.
[void]$ws.Range("A10").Select()
